# modify System power consumption
# Insert a new "PK" worksheet between "Topology" and "10.4.L.5.2", fill it
# with the notes/benchmark-comparison content, and move the active/selected
# tab to the last sheet ("40.10.L.5.2").

$wb = $excel.ActiveWorkbook

# --- Insert the new "PK" sheet right before "10.4.L.5.2" -------------------
$before = $wb.Worksheets.Item("10.4.L.5.2")
$pk = $wb.Worksheets.Add($before)
$pk.Name = "PK"

# --- Fill the PK sheet content (row order matches original authoring order)
$pk.Range("A2").Value = "不知道preferrence vector"
$pk.Range("B2").Value = "知道preferrence vector"
$pk.Range("A3").Value = "不須決定最大連線數"
$pk.Range("B3").Value = "需要依照L計算"
$pk.Range("A4").Value = "以最大EE為目標"
$pk.Range("B4").Value = "只有caching是最小化Hit Rate來最小P_sys，間接最小化EE`n畫出所有L，找一個最好的EE跟RL PK"
$pk.Range("B5").Value = "可以嘗試SINR Benchmark"
$pk.Range("C6").Value = "坤霖: Req的itr變大, 簡老師: Req先不變, implement time variant channel "
$pk.Range("A1").Value = "RL"
$pk.Range("B1").Value = "Benchmark(知道越多越好)"

# B4 holds the long paragraph -> wrap text + an explicit row height like the
# authored sheet.
$pk.Range("B4").WrapText = $true
$pk.Rows(4).RowHeight = 63

# Column sizing to roughly match the authored "best fit" widths.
$pk.Columns(1).ColumnWidth = 23.7109375
$pk.Columns(2).ColumnWidth = 26.42578125

# Page setup (A4/portrait single page) like the other data sheets.
$pk.PageSetup.PaperSize = 9
$pk.PageSetup.Orientation = 1

$pk.Range("B12").Select()

# --- Move the active tab to the last sheet ---------------------------------
$wb.Worksheets.Item("40.10.L.5.2").Activate()
